$p = $ppt.ActivePresentation
$layouts = $p.SlideMaster.CustomLayouts
$l1 = $layouts.Item(1)
$sh = $l1.Shapes.Item(4)
Write-Output $sh.TextFrame.TextRange.Text
$sh.TextFrame.TextRange.Text = "8/26/19"
Write-Output $sh.TextFrame.TextRange.Text
